# Generate Report for Handoff
# -----------------------------------------------------------------------
# The localization-status report gains two new tracked files
# (198d1cd1-...png and b0c94c0e-...png) that sit around the existing
# 63351912-...md entry (previously named ab744e59-...md). Each of the
# three worksheets (Overview, zh-cn, de-de) grows from 3 data rows to 5.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- commonly reused literal strings -----------------------------------
$readyForHandoff   = "Ready for handoff"
$notLocalized      = "Not to be localized"
$configName        = ".localization-config"

$pngA              = "198d1cd1-7709-4945-9f0c-9710a0857768.png"
$mdName            = "63351912-1311-4686-957b-afd972ca2238.md"
$pngB              = "b0c94c0e-346c-479c-99e7-8a38b61a387a.png"

$zhPngXlf          = "0dd162be70daa77935309e9e4ac503b0b59ebb1b.png"
$zhMdXlf           = "63351912-1311-4686-957b-afd972ca2238.d5be2bd33796896dd0992e34484c1a0b3dd7dc38.zh-cn.xlf"
$dePngXlf          = "0dd162be70daa77935309e9e4ac503b0b59ebb1b.png"
$deMdXlf           = "63351912-1311-4686-957b-afd972ca2238.d5be2bd33796896dd0992e34484c1a0b3dd7dc38.de-de.xlf"
$pngBXlf           = "20f61260df2cbce151f23658049786d0c743414b.png"

$zhTime            = "2016-03-10 21:09:04"
$deTime            = "2016-03-10 21:09:10"
$zeroTime          = "0001-01-01 00:00:00"

$includeTxt        = "Include"
$ignoredTxt        = "Ignored"
$isDependencyTxt   = "IsDependency"
$dependencyFromMd  = "e2e\63351912-1311-4686-957b-afd972ca2238.md"

# Base (placeholder, but pattern-consistent) hyperlink target prefixes.
$rawCommit   = "a1b2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7f8a9b0"
$srcBase     = "https://github.com/OpenLocalizationTest/oltest/blob/$rawCommit/e2e"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/68dd313bc44326ec1503c2f506e3655bb2af9dd1/.localization-config"
$zhHandoffCommit = "d5be2bd33796896dd0992e34484c1a0b3dd7dc38"
$deHandoffCommit = "d5be2bd33796896dd0992e34484c1a0b3dd7dc38"
$zhHtBase    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHtBase    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

function Set-RowFromArray($ws, [int]$row, [string[]]$values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        if ($null -ne $values[$i]) {
            $ws.Cells.Item($row, $i + 1).Value = $values[$i]
        }
    }
}

# Matches the workbook's existing "HyperLink" cell style (underlined,
# cornflower blue FF6495ED) so newly-added links look like the originals
# instead of picking up the generic theme-colored auto style.
function Set-HyperlinkLook($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

function Add-Hyperlink($ws, [string]$cellRef, [string]$address, [string]$display) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $display) | Out-Null
    Set-HyperlinkLook $ws.Range($cellRef)
}

# =========================================================================
# Sheet "Overview" (A:C) -- File Name / zh-cn / de-de
# =========================================================================
$wsOv = $wb.Worksheets.Item("Overview")

# Drop existing hyperlinks before shuffling rows around -- the engine does
# not re-anchor hyperlink ranges when rows are inserted/shifted, so we
# remove them first and re-create them at their final location afterwards.
$wsOv.Hyperlinks.Delete()

# Insert two fresh rows ahead of the ".localization-config" row (old row 3)
# so it lands on row 5 in the final layout.
$wsOv.Rows("3:4").Insert()

Set-RowFromArray $wsOv 2 @($pngA, $readyForHandoff, $readyForHandoff)
Set-RowFromArray $wsOv 3 @($mdName, $readyForHandoff, $readyForHandoff)
Set-RowFromArray $wsOv 4 @($pngB, $readyForHandoff, $readyForHandoff)
Set-RowFromArray $wsOv 5 @($configName, $notLocalized, $notLocalized)

Add-Hyperlink $wsOv "A2" "$srcBase/$pngA"  $pngA
Add-Hyperlink $wsOv "A3" "$srcBase/$mdName" $mdName
Add-Hyperlink $wsOv "A4" "$srcBase/$pngB"  $pngB
Add-Hyperlink $wsOv "A5" $cfgUrl           $configName

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()
$wsZh.Rows("3:4").Insert()

# Row 2 -- existing ab744e59 entry becomes the 198d1cd1 png entry.
Set-RowFromArray $wsZh 2 @($pngA, $readyForHandoff, $zhPngXlf, $zhTime, $null, $null, $zeroTime, $isDependencyTxt, $dependencyFromMd)
# Row 3 (new) -- the .md entry, now with its own handoff xlf.
Set-RowFromArray $wsZh 3 @($mdName, $readyForHandoff, $zhMdXlf, $zhTime, $null, $null, $zeroTime, $includeTxt, $null)
# Row 4 (new) -- the b0c94c0e png entry.
Set-RowFromArray $wsZh 4 @($pngB, $readyForHandoff, $pngBXlf, $zhTime, $null, $null, $zeroTime, $isDependencyTxt, $dependencyFromMd)
# Row 5 -- shifted-down .localization-config row (unchanged content).
Set-RowFromArray $wsZh 5 @($configName, $notLocalized, $null, $zeroTime, $null, $null, $zeroTime, $ignoredTxt, $null)

Add-Hyperlink $wsZh "A2" "$srcBase/$pngA"      $pngA
Add-Hyperlink $wsZh "C2" "$zhHtBase/$zhPngXlf" $zhPngXlf
Add-Hyperlink $wsZh "A3" "$srcBase/$mdName"    $mdName
Add-Hyperlink $wsZh "C3" "$zhHtBase/$zhMdXlf"  $zhMdXlf
Add-Hyperlink $wsZh "A4" "$srcBase/$pngB"      $pngB
Add-Hyperlink $wsZh "C4" "$zhHtBase/$pngBXlf"  $pngBXlf
Add-Hyperlink $wsZh "A5" $cfgUrl               $configName

# =========================================================================
# Sheet "de-de"
# =========================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()
$wsDe.Rows("3:4").Insert()

Set-RowFromArray $wsDe 2 @($pngA, $readyForHandoff, $dePngXlf, $deTime, $null, $null, $zeroTime, $isDependencyTxt, $dependencyFromMd)
Set-RowFromArray $wsDe 3 @($mdName, $readyForHandoff, $deMdXlf, $deTime, $null, $null, $zeroTime, $includeTxt, $null)
Set-RowFromArray $wsDe 4 @($pngB, $readyForHandoff, $pngBXlf, $deTime, $null, $null, $zeroTime, $isDependencyTxt, $dependencyFromMd)
Set-RowFromArray $wsDe 5 @($configName, $notLocalized, $null, $zeroTime, $null, $null, $zeroTime, $ignoredTxt, $null)

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$srcBase/$pngA", "", "", $pngA)                 | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "$deHtBase/$dePngXlf", "", "", $dePngXlf)         | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$srcBase/$mdName", "", "", $mdName)              | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "$deHtBase/$deMdXlf", "", "", $deMdXlf)           | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$srcBase/$pngB", "", "", $pngB)                  | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "$deHtBase/$pngBXlf", "", "", $pngBXlf)           | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $cfgUrl, "", "", $configName)                     | Out-Null
